# "Generate Report for Handback" - populate the Latest Target File / Latest
# Handback File columns (F/G) for each localized-language sheet, refresh the
# Status + Latest Handback DateTime columns to reflect a completed handback.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$zhHandbackDate = "2016-03-12 04:37:47"
$deRow2HandbackDate = "2016-03-12 04:37:47"
$deRow3HandbackDate = "2016-03-12 04:37:52"

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/ab51247e39af4e9dde8ce4b242cb91cfa8039f8e/e2e"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2 (5bbc677a...)
$ws.Range("C2").Value = $statusHandedBack
$ws.Hyperlinks.Add($ws.Range("F2"), "$repoBase/5bbc677a-3fb8-45cd-aabc-4770dab871d6.md", [Type]::Missing, [Type]::Missing, "5bbc677a-3fb8-45cd-aabc-4770dab871d6.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/acdf47183e4bdc1314197d481c330fb1efdb11a6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5bbc677a-3fb8-45cd-aabc-4770dab871d6.1ef34bcc4ae7ace7bbd6e11cc9425393498071a5.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "5bbc677a-3fb8-45cd-aabc-4770dab871d6.1ef34bcc4ae7ace7bbd6e11cc9425393498071a5.zh-cn.xlf")
$ws.Range("H2").Value = $zhHandbackDate

# Row 3 (c8d950fd...)
$ws.Range("C3").Value = $statusHandedBack
$ws.Hyperlinks.Add($ws.Range("F3"), "$repoBase/c8d950fd-4bee-4b9b-bb81-55f430cef6fc.md", [Type]::Missing, [Type]::Missing, "c8d950fd-4bee-4b9b-bb81-55f430cef6fc.md")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/acdf47183e4bdc1314197d481c330fb1efdb11a6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c8d950fd-4bee-4b9b-bb81-55f430cef6fc.1acd3de55168fb24aefb4ca660fa240cebebd17a.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "c8d950fd-4bee-4b9b-bb81-55f430cef6fc.1acd3de55168fb24aefb4ca660fa240cebebd17a.zh-cn.xlf")
$ws.Range("H3").Value = $zhHandbackDate

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 2 (5bbc677a...)
$ws.Range("C2").Value = $statusHandedBack
$ws.Hyperlinks.Add($ws.Range("F2"), "$repoBase/5bbc677a-3fb8-45cd-aabc-4770dab871d6.md", [Type]::Missing, [Type]::Missing, "5bbc677a-3fb8-45cd-aabc-4770dab871d6.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ab1c2f0615a4bfcba483d5d799263614252538e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5bbc677a-3fb8-45cd-aabc-4770dab871d6.1ef34bcc4ae7ace7bbd6e11cc9425393498071a5.de-de.xlf", [Type]::Missing, [Type]::Missing, "5bbc677a-3fb8-45cd-aabc-4770dab871d6.1ef34bcc4ae7ace7bbd6e11cc9425393498071a5.de-de.xlf")
$ws.Range("H2").Value = $deRow2HandbackDate

# Row 3 (c8d950fd...)
$ws.Range("C3").Value = $statusHandedBack
$ws.Hyperlinks.Add($ws.Range("F3"), "$repoBase/c8d950fd-4bee-4b9b-bb81-55f430cef6fc.md", [Type]::Missing, [Type]::Missing, "c8d950fd-4bee-4b9b-bb81-55f430cef6fc.md")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ab1c2f0615a4bfcba483d5d799263614252538e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c8d950fd-4bee-4b9b-bb81-55f430cef6fc.1acd3de55168fb24aefb4ca660fa240cebebd17a.de-de.xlf", [Type]::Missing, [Type]::Missing, "c8d950fd-4bee-4b9b-bb81-55f430cef6fc.1acd3de55168fb24aefb4ca660fa240cebebd17a.de-de.xlf")
$ws.Range("H3").Value = $deRow3HandbackDate
